$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.571.25'
$ws.Range('D3').Value = '1.884.44'
$ws.Range('E3').Value = '  +1.46%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.47'
$ws.Range('E5').Value = '  +6.07%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4775'
$ws.Range('E7').Value = '  +2.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2925'
$ws.Range('E8').Value = '  +3.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06537'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.10'
$ws.Range('E10').Value = '  +5.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07723'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.91'
$ws.Range('E12').Value = '  +4.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7434'
$ws.Range('E13').Value = '  +9.31%  '
$ws.Range('D14').Value = '1.880.67'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.160'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '275.17'
$ws.Range('E16').Value = '  +3.65%  '
$ws.Range('D17').Value = '30.566.43'
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.51'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007580'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = '2.130.51'
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.273'
$ws.Range('E23').Value = '  +2.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.202'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.350'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.43'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('E27').Value = '  +2.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.952'
$ws.Range('E28').Value = '  +3.64%  '
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09989'
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('E31').Value = '  +4.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.334'
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.129'
$ws.Range('E33').Value = '  +3.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04813'
$ws.Range('E34').Value = '  +3.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.133'
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7038'
$ws.Range('E36').Value = '  +2.82%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01876'
$ws.Range('E38').Value = '  +3.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.752'
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.341'
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.985'
$ws.Range('E41').Value = '  +5.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '71.66'
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4240'
$ws.Range('E43').Value = '  +4.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8421'
$ws.Range('E44').Value = '  +1.31%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.85'
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.356'
$ws.Range('E47').Value = '  +2.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.115'
$ws.Range('E48').Value = '  +2.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.74'
$ws.Range('E49').Value = '  +4.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '919.30'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3921'
$ws.Range('E51').Value = '  +4.94%  '
